$wb = $excel.ActiveWorkbook

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 990.0909
$ws.Range("I111").Value = 979.1
$ws.Range("K111").Value = 2937.3
$ws.Range("M111").Value = 129.6999999999998

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 13159561
$ws.Range("I112").Value = 2816.6667
$ws.Range("J112").Value = 19231904
$ws.Range("K112").Value = 8450.000100000001
$ws.Range("L112").Value = 57695712
$ws.Range("M112").Value = -7342.000100000001
$ws.Range("N112").Value = -57697928

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4678.8335
$ws.Range("I113").Value = 4457.143
$ws.Range("J113").Value = 4989.2
$ws.Range("K113").Value = 4457.143
$ws.Range("L113").Value = 4989.2
$ws.Range("M113").Value = -1203.143
$ws.Range("N113").Value = -11497.2

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5729.75
$ws.Range("I137").Value = 7167.7144
$ws.Range("J137").Value = 3716.6
$ws.Range("K137").Value = 21503.1432
$ws.Range("L137").Value = 11149.8
$ws.Range("M137").Value = -18953.1432
$ws.Range("N137").Value = -16249.8

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3512.7092
$ws.Range("I138").Value = 1864.9807
$ws.Range("J138").Value = 6032.7646
$ws.Range("K138").Value = 5594.9421
$ws.Range("L138").Value = 18098.2938
$ws.Range("M138").Value = -454.9421000000002
$ws.Range("N138").Value = -28378.2938

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 478106.8
$ws.Range("I141").Value = 1248.5555
$ws.Range("J141").Value = 1336451.8
$ws.Range("K141").Value = 3745.6665
$ws.Range("L141").Value = 4009355.4
$ws.Range("M141").Value = 1434.3335
$ws.Range("N141").Value = -4019715.4

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14680.9
$ws.Range("I32").Value = 12472.147
$ws.Range("K32").Value = 12472.147
$ws.Range("M32").Value = -12185.147

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1781.6792
$ws.Range("I122").Value = 1389.6511
$ws.Range("J122").Value = 3467.4
$ws.Range("K122").Value = 4168.9533
$ws.Range("L122").Value = 10402.2
$ws.Range("M122").Value = -1718.9533
$ws.Range("N122").Value = -15302.2

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2902.4285
$ws.Range("I134").Value = 2860.2285
$ws.Range("J134").Value = 3113.4285
$ws.Range("K134").Value = 8580.6855
$ws.Range("L134").Value = 9340.2855
$ws.Range("M134").Value = -6045.6855
$ws.Range("N134").Value = -14410.2855

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 730001.4399999999
$ws.Range("I6").Value = 2000000.4
$ws.Range("K6").Value = 2000000.4
$ws.Range("M6").Value = -1999887.4

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6968.9565
$ws.Range("I31").Value = 8978
$ws.Range("J31").Value = 6546
$ws.Range("K31").Value = 8978
$ws.Range("L31").Value = 6546
$ws.Range("M31").Value = -8683
$ws.Range("N31").Value = -7136

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6968.9565
$ws.Range("I34").Value = 8978
$ws.Range("J34").Value = 6546
$ws.Range("K34").Value = 8978
$ws.Range("L34").Value = 6546
$ws.Range("M34").Value = -8776
$ws.Range("N34").Value = -6950

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8067417
$ws.Range("I58").Value = 1403.9474
$ws.Range("K58").Value = 1403.9474
$ws.Range("M58").Value = -1200.9474

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3199.3845
$ws.Range("I86").Value = 1949.5
$ws.Range("J86").Value = 4270.7144
$ws.Range("K86").Value = 1949.5
$ws.Range("L86").Value = 4270.7144
$ws.Range("M86").Value = -826.5
$ws.Range("N86").Value = -6516.7144

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3199.3845
$ws.Range("I89").Value = 1949.5
$ws.Range("J89").Value = 4270.7144
$ws.Range("K89").Value = 9747.5
$ws.Range("L89").Value = 21353.572
$ws.Range("M89").Value = -4131.5
$ws.Range("N89").Value = -32585.572

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5067.6665
$ws.Range("I99").Value = 2775
$ws.Range("J99").Value = 6214
$ws.Range("K99").Value = 2775
$ws.Range("L99").Value = 6214
$ws.Range("M99").Value = -1277
$ws.Range("N99").Value = -9210

# CRP row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -24920

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5067.6665
$ws.Range("I126").Value = 2775
$ws.Range("J126").Value = 6214
$ws.Range("K126").Value = 8325
$ws.Range("L126").Value = 18642
$ws.Range("M126").Value = -5855
$ws.Range("N126").Value = -23582

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8067417
$ws.Range("I136").Value = 1403.9474
$ws.Range("K136").Value = 4211.8422
$ws.Range("M136").Value = -1661.8422

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 31051.725
$ws.Range("J141").Value = 31051.725
$ws.Range("L141").Value = 31051.725
$ws.Range("N141").Value = -41411.725

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 610.6
$ws.Range("I23").Value = 101
$ws.Range("J23").Value = 667.2222
$ws.Range("K23").Value = 303
$ws.Range("L23").Value = 2001.6666
$ws.Range("M23").Value = -68
$ws.Range("N23").Value = -2471.6666

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5659.0835
$ws.Range("I56").Value = 5659.0835
$ws.Range("K56").Value = 5659.0835
$ws.Range("M56").Value = -5129.0835

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 833.51514
$ws.Range("I113").Value = 634
$ws.Range("J113").Value = 869.1429000000001
$ws.Range("K113").Value = 1902
$ws.Range("L113").Value = 2607.4287
$ws.Range("M113").Value = 268
$ws.Range("N113").Value = -6947.4287

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 9891.151
$ws.Range("I121").Value = 218.75
$ws.Range("J121").Value = 18994.588
$ws.Range("K121").Value = 656.25
$ws.Range("L121").Value = 56983.764
$ws.Range("M121").Value = 653.75
$ws.Range("N121").Value = -59603.764

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1219.6364
$ws.Range("I122").Value = 552.125
$ws.Range("K122").Value = 4969.125
$ws.Range("M122").Value = -2519.125

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3252.1304
$ws.Range("J137").Value = 4495.5713
$ws.Range("L137").Value = 13486.7139
$ws.Range("M137").Value = -3024.375
$ws.Range("N137").Value = -23686.7139

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3150
$ws.Range("I80").Value = 3150
$ws.Range("K80").Value = 3150
$ws.Range("M80").Value = -2152

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3150
$ws.Range("I83").Value = 3150
$ws.Range("K83").Value = 15750
$ws.Range("M83").Value = -10758

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2751
$ws.Range("I82").Value = 1676.6666
$ws.Range("J82").Value = 4362.5
$ws.Range("K82").Value = 1676.6666
$ws.Range("L82").Value = 4362.5
$ws.Range("M82").Value = -1315.6666
$ws.Range("N82").Value = -5084.5

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2751
$ws.Range("I85").Value = 1676.6666
$ws.Range("J85").Value = 4362.5
$ws.Range("K85").Value = 1676.6666
$ws.Range("L85").Value = 4362.5
$ws.Range("M85").Value = -428.6666
$ws.Range("N85").Value = -6858.5

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2348.5715
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 2875.5557
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 2875.5557
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -3957.5557

# LTW row 101
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 19693.666
$ws.Range("J101").Value = 19693.666
$ws.Range("L101").Value = 19693.666
$ws.Range("N101").Value = -26183.666

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1831.25
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 35635.773
$ws.Range("I126").Value = 45774.87
$ws.Range("J126").Value = 6485.875
$ws.Range("K126").Value = 137324.61
$ws.Range("L126").Value = 19457.625
$ws.Range("M126").Value = -134854.61
$ws.Range("N126").Value = -24397.625

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2374.6667
$ws.Range("I136").Value = 1918.8214
$ws.Range("J136").Value = 3286.3572
$ws.Range("K136").Value = 5756.4642
$ws.Range("L136").Value = 9859.071599999999
$ws.Range("M136").Value = -3206.4642
$ws.Range("N136").Value = -14959.0716
